$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data per the diff (Sun Mar 5 2023 GitHub Actions run)
# Force text format on target cells so numeric-looking strings (e.g. "1.003") are preserved as text

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.444.05'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.25%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.572.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.28'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3743'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.88'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3403'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.147'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07576'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.85%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.37'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.989'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.954'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.573.19'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001124'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.12'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.85%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06732'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.35%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.276'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.42'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.21'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.93%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.451.67'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.319'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.595'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.94%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.16'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '148.68'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.60%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.002'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.03'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.748.70'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.137'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.979'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.866'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.63%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08435'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.385'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02465'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.44%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2298'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06539'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.64%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.492'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.39'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6290'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.16%  '

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.13'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.09%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.814'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5873'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.096'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.10'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.62%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.230'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.74%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07332'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.02%  '
